# Applies the "demo dataset for the NoiseCompressor project" update:
# the stem_infos sheet is populated with the full 30-stem dataset
# (replacing the 2-row sample that only contained stems 20709 and 22729).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stem_infos")

# --- Prepare cell formatting for the rows that will receive data -------
# Column A uses two different styles in the final sheet:
#  - rows 2-18 (ids 10001-14648) use the same style as the header cell A1
#  - rows 19-31 (ids 20506-23436) use the same style as the pre-existing
#    data cell A2 (the style historically reserved for the 20000+ stems)
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A19:A31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A1").Copy() | Out-Null
$ws.Range("A2:A18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Columns B, C and D keep the same formatting as the existing data rows
# for every new row (rows 4-31 don't have it yet).
$ws.Range("B2:D2").Copy() | Out-Null
$ws.Range("B4:D31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Fill in the full dataset ------------------------------------------
$rows = @(
    @{Row=2; A=10001; B="10001.ply"; C=0.0101; D=1},
    @{Row=3; A=10002; B="10002.ply"; C=0.0226; D=1},
    @{Row=4; A=10004; B="10004.ply"; C=0.0463; D=1},
    @{Row=5; A=10006; B="10006.ply"; C=0.06570000000000001; D=1},
    @{Row=6; A=10008; B="10008.ply"; C=0.0767; D=1},
    @{Row=7; A=10010; B="10010.ply"; C=0.0867; D=1},
    @{Row=8; A=10809; B="10809.ply"; C=0.085; D=1},
    @{Row=9; A=11011; B="11011.ply"; C=0.105; D=1},
    @{Row=10; A=11315; B="11315.ply"; C=0.14; D=1},
    @{Row=11; A=11820; B="11820.ply"; C=0.19; D=1},
    @{Row=12; A=12323; B="12323.ply"; C=0.23; D=1},
    @{Row=13; A=12628; B="12628.ply"; C=0.27; D=1},
    @{Row=14; A=13031; B="13031.ply"; C=0.305; D=1},
    @{Row=15; A=13335; B="13335.ply"; C=0.34; D=1},
    @{Row=16; A=13638; B="13638.ply"; C=0.37; D=1},
    @{Row=17; A=14445; B="14445.ply"; C=0.445; D=1},
    @{Row=18; A=14648; B="14648.ply"; C=0.47; D=1},
    @{Row=19; A=20506; B="20506.ply"; C=0.055; D=1},
    @{Row=20; A=20607; B="20607.ply"; C=0.065; D=1},
    @{Row=21; A=20709; B="20709.ply"; C=0.08; D=1},
    @{Row=22; A=20809; B="20809.ply"; C=0.085; D=1},
    @{Row=23; A=21112; B="21112.ply"; C=0.115; D=1},
    @{Row=24; A=21315; B="21315.ply"; C=0.14; D=1},
    @{Row=25; A=21617; B="21617.ply"; C=0.165; D=1},
    @{Row=26; A=21719; B="21719.ply"; C=0.18; D=1},
    @{Row=27; A=22526; B="22526.ply"; C=0.255; D=1},
    @{Row=28; A=22729; B="22729.ply"; C=0.28; D=1},
    @{Row=29; A=22830; B="22830.ply"; C=0.29; D=1},
    @{Row=30; A=23132; B="23132.ply"; C=0.315; D=1},
    @{Row=31; A=23436; B="23436.ply"; C=0.35; D=1}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}

$ws.Range("G11").Select() | Out-Null
